# Apply cryptocurrency price/volume updates as described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "80.731.65"
$ws.Range("E2").Value = "  +5.50%  "

# Row 3
$ws.Range("D3").Value = "3.229.72"
$ws.Range("E3").Value = "  +5.14%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.12%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.50"
$ws.Range("E5").Value = "  +7.00%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "637.26"
$ws.Range("E6").Value = "  +3.14%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.284"
$ws.Range("E7").Value = "  +35.67%  "

# Row 8
$ws.Range("E8").Value = "  -0.12%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.611"
$ws.Range("E9").Value = "  +10.96%  "

# Row 10
$ws.Range("D10").Value = "3.227.50"
$ws.Range("E10").Value = "  +5.16%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.618"
$ws.Range("E11").Value = "  +40.39%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000266"
$ws.Range("E12").Value = "  +37.64%  "

# Row 13
$ws.Range("E13").Value = "  +3.64%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.46"
$ws.Range("E14").Value = "  +4.18%  "

# Row 15
$ws.Range("D15").Value = "3.820.76"
$ws.Range("E15").Value = "  +5.85%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "32.90"
$ws.Range("E16").Value = "  +13.40%  "

# Row 17
$ws.Range("D17").Value = "80.641.01"
$ws.Range("E17").Value = "  +5.46%  "

# Row 18
$ws.Range("D18").Value = "3.226.35"
$ws.Range("E18").Value = "  +5.38%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "14.73"
$ws.Range("E19").Value = "  +8.47%  "

# Row 20
$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "9.48"
$ws.Range("E20").Value = "  +6.19%  "

# Row 21
$ws.Range("B21").Value = "BitcoinCash"
$ws.Range("C21").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "451.40"
$ws.Range("E21").Value = "  +18.26%  "

# Row 22
$ws.Range("B22").Value = "SuiNetwork"
$ws.Range("C22").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.02"
$ws.Range("E22").Value = "  +22.58%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.39"
$ws.Range("E23").Value = "  +22.28%  "

# Row 24
$ws.Range("B24").Value = "WrappedeETH"
$ws.Range("C24").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D24").Value = "3.393.86"
$ws.Range("E24").Value = "  +5.60%  "

# Row 25
$ws.Range("B25").Value = "NEARProtocol"
$ws.Range("C25").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.83"
$ws.Range("E25").Value = "  +11.17%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "78.07"
$ws.Range("E26").Value = "  +7.63%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.12"
$ws.Range("E27").Value = "  +12.75%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0000128"
$ws.Range("E28").Value = "  +18.24%  "

# Row 29
$ws.Range("E29").Value = "  +0.03%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.35"
$ws.Range("E30").Value = "  +12.31%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.996"
$ws.Range("E31").Value = "  -0.71%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "567.00"
$ws.Range("E32").Value = "  +13.70%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.52"
$ws.Range("E33").Value = "  +8.98%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.156"
$ws.Range("E34").Value = "  +25.26%  "

# Row 35
$ws.Range("E35").Value = "  +7.33%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "24.04"
$ws.Range("E36").Value = "  +16.03%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.125"
$ws.Range("E37").Value = "  +20.81%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.423"
$ws.Range("E38").Value = "  +11.74%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.999"
$ws.Range("E39").Value = "  -0.10%  "

# Row 40
$ws.Range("B40").Value = "Monero"
$ws.Range("C40").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "165.93"
$ws.Range("E40").Value = "  +2.05%  "

# Row 41
$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.88"
$ws.Range("E41").Value = "  +14.58%  "

# Row 42
$ws.Range("E42").Value = "  +1.54%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "193.67"
$ws.Range("E43").Value = "  +0.36%  "

# Row 44
$ws.Range("E44").Value = "  +0.02%  "

# Row 45
$ws.Range("B45").Value = "Stacks"
$ws.Range("C45").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.86"
$ws.Range("E45").Value = "  +13.34%  "

# Row 46
$ws.Range("B46").Value = "dogwifhat"
$ws.Range("C46").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.78"
$ws.Range("E46").Value = "  +14.17%  "

# Row 47
$ws.Range("E47").Value = "  +8.88%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.805"
$ws.Range("E48").Value = "  +2.12%  "

# Row 49
$ws.Range("B49").Value = "OKB"
$ws.Range("C49").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "43.84"
$ws.Range("E49").Value = "  +6.31%  "

# Row 50
$ws.Range("B50").Value = "Filecoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "4.39"
$ws.Range("E50").Value = "  +13.42%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.649"
$ws.Range("E51").Value = "  +8.81%  "

